$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 9,20
$data[0,0] = "ECs"
$data[0,1] = "Nlgn3"
$data[0,2] = "Nrxn1"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 0.2210183333333333
$data[0,7] = 0.6630550000000001
$data[0,8] = 0.111623447733668
$data[0,9] = 0.111623447733668
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.002858666666666667
$data[0,13] = 0.008576
$data[0,14] = 0.002669495535069502
$data[0,15] = 0.002669495535069501
$data[0,16] = 0.0006318177422222223
$data[0,17] = 0.005686359680000001
$data[0,18] = 0.0002979782953340906
$data[0,19] = 0.0002979782953340906
$data[1,0] = "ECs"
$data[1,1] = "Nlgn3"
$data[1,2] = "Nrxn1"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 0.2210183333333333
$data[1,7] = 0.6630550000000001
$data[1,8] = 0.111623447733668
$data[1,9] = 0.111623447733668
$data[1,10] = 2
$data[1,11] = 0.6666666666666666
$data[1,12] = 0.07823633333333334
$data[1,13] = 0.234709
$data[1,14] = 0.07305907503971872
$data[1,15] = 0.0730590750397187
$data[1,16] = 0.01729166399944445
$data[1,17] = 0.155624975995
$data[1,18] = 0.008155105844166172
$data[1,19] = 0.00815510584416617
$data[2,0] = "ECs"
$data[2,1] = "Nlgn3"
$data[2,2] = "Nrxn1"
$data[2,3] = "MuSCs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 0.2210183333333333
$data[2,7] = 0.6630550000000001
$data[2,8] = 0.111623447733668
$data[2,9] = 0.111623447733668
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.989769
$data[2,13] = 2.969307
$data[2,14] = 0.9242714294252118
$data[2,15] = 0.9242714294252118
$data[2,16] = 0.218757094765
$data[2,17] = 1.968813852885
$data[2,18] = 0.1031703635941678
$data[2,19] = 0.1031703635941678
$data[3,0] = "FAPs"
$data[3,1] = "Nlgn3"
$data[3,2] = "Nrxn1"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 0.2210183333333333
$data[3,7] = 0.6630550000000001
$data[3,8] = 0.5186760166697389
$data[3,9] = 0.5186760166697389
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.002858666666666667
$data[3,13] = 0.008576
$data[3,14] = 0.002669495535069502
$data[3,15] = 0.002669495535069501
$data[3,16] = 0.002935841137777778
$data[3,17] = 0.02642257024
$data[3,18] = 0.001384603310647503
$data[3,19] = 0.001384603310647502
$data[4,0] = "FAPs"
$data[4,1] = "Nlgn3"
$data[4,2] = "Nrxn1"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 1.026996666666667
$data[4,7] = 3.08099
$data[4,8] = 0.5186760166697389
$data[4,9] = 0.5186760166697389
$data[4,10] = 2
$data[4,11] = 0.6666666666666666
$data[4,12] = 0.07823633333333334
$data[4,13] = 0.234709
$data[4,14] = 0.07305907503971872
$data[4,15] = 0.0730590750397187
$data[4,16] = 0.08034845354555555
$data[4,17] = 0.72313608191
$data[4,18] = 0.03789399002317685
$data[4,19] = 0.03789399002317684
$data[5,0] = "FAPs"
$data[5,1] = "Nlgn3"
$data[5,2] = "Nrxn1"
$data[5,3] = "MuSCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 1.026996666666667
$data[5,7] = 3.08099
$data[5,8] = 0.5186760166697389
$data[5,9] = 0.5186760166697389
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.989769
$data[5,13] = 2.969307
$data[5,14] = 0.9242714294252118
$data[5,15] = 0.9242714294252118
$data[5,16] = 1.01648946377
$data[5,17] = 9.14840517393
$data[5,18] = 0.4793974233359146
$data[5,19] = 0.4793974233359146
$data[6,0] = "MuSCs"
$data[6,1] = "Nlgn3"
$data[6,2] = "Nrxn1"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 0.73202
$data[6,7] = 2.19606
$data[6,8] = 0.369700535596593
$data[6,9] = 0.369700535596593
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.002858666666666667
$data[6,13] = 0.008576
$data[6,14] = 0.002669495535069502
$data[6,15] = 0.002669495535069501
$data[6,16] = 0.002092601173333334
$data[6,17] = 0.01883341056
$data[6,18] = 0.0009869139290879085
$data[6,19] = 0.0009869139290879082
$data[7,0] = "MuSCs"
$data[7,1] = "Nlgn3"
$data[7,2] = "Nrxn1"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.73202
$data[7,7] = 2.19606
$data[7,8] = 0.369700535596593
$data[7,9] = 0.369700535596593
$data[7,10] = 2
$data[7,11] = 0.6666666666666666
$data[7,12] = 0.07823633333333334
$data[7,13] = 0.234709
$data[7,14] = 0.07305907503971872
$data[7,15] = 0.0730590750397187
$data[7,16] = 0.05727056072666667
$data[7,17] = 0.5154350465400001
$data[7,18] = 0.02700997917237569
$data[7,19] = 0.02700997917237568
$data[8,0] = "MuSCs"
$data[8,1] = "Nlgn3"
$data[8,2] = "Nrxn1"
$data[8,3] = "MuSCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 0.73202
$data[8,7] = 2.19606
$data[8,8] = 0.369700535596593
$data[8,9] = 0.369700535596593
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.989769
$data[8,13] = 2.969307
$data[8,14] = 0.9242714294252118
$data[8,15] = 0.9242714294252118
$data[8,16] = 0.72453070338
$data[8,17] = 6.52077633042
$data[8,18] = 0.3417036424951294
$data[8,19] = 0.3417036424951294

$ws.Range("A2:T10").Value = $data

